$p = $ppt.ActivePresentation

# --- Slide 3: reposition "Content Placeholder 2" textbox ---
# EMU (-69960,1700262) -> (0,1012105); Shape.Left/Top are in points (1pt = 12700 EMU)
$s3 = $p.Slides.Item(3)
$contentBox = $s3.Shapes.Item(2)
$contentBox.Left = 0.000000000000
$contentBox.Top = 79.693346436005

# --- Slide 5: reposition the three ranking pictures ---
$s5 = $p.Slides.Item(5)

$pic6 = $s5.Shapes.Item(2)
$pic6.Left = -0.000078840156
$pic6.Top = 97.930984497070

$pic7 = $s5.Shapes.Item(3)
$pic7.Left = 325.745315551758
$pic7.Top = 102.131298065186

$pic8 = $s5.Shapes.Item(4)
$pic8.Left = 648.573028564453
$pic8.Top = 102.131298065186

# --- Slide 5: remove the two curved-down arrow callouts ---
# Delete from the end so indices of the remaining shapes stay valid.
$s5.Shapes.Item(6).Delete()
$s5.Shapes.Item(5).Delete()
